$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Info" column header (E1): bold, centered/top aligned like the other
# headers, but with just a left border, plus a wider column.
$ws.Range("E1").Value = "Info"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").VerticalAlignment = -4160
$ws.Range("E1").Borders.Item(7).LineStyle = 1
$ws.Range("E1").ColumnWidth = 45

# "Info" values for the two pre-existing rows (2-3) - their A/B/C values are
# untouched/unchanged.
$ws.Cells.Item(2, 5).Value = "Blazer 9mm 124"
$ws.Cells.Item(3, 5).Value = "PMC BRONZE 223"

# New rows of ammo data (rows 4-11).
$rows = @(
    @{ Row = 4;  A = 3877002037092;  B = 120;  C = 79.99;  E = "IGMAN 7.62X39" },
    @{ Row = 5;  A = 29465025922;    B = 250;  C = 111.99; E = "FEDERAL 12GA" },
    @{ Row = 6;  A = 47700346007;    B = 250;  C = 111.99; E = "FEDERAL 12GA" },
    @{ Row = 7;  A = 20892225008;    B = 150;  C = 79.99;  E = "WINCHESTER 5.56 150RD " },
    @{ Row = 8;  A = 604544634211;   B = 200;  C = 65.99;  E = "FEDERAL 9MM 200RD" },
    @{ Row = 9;  A = 604544695861;   B = 1100; C = 79.99;  E = "FEDERAL 22LR BLACK PACK" },
    @{ Row = 10; A = 754908200313;   B = 1000; C = 499.99; E = "CBC MAGTECH 5.56" },
    @{ Row = 11; A = 76683052308;    B = 1000; C = 429.99; E = "BLAZER 45 ACP" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 5).Value = $r.E
}

# Rows 4-11 need the same "UPC" column styling (greyed Arial font, via the
# style already used by A2:A3) applied to column A.
$ws.Range("A2").Copy()
$ws.Range("A4:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E11").Select()
